$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.327.71"
$ws.Range("E2").Value = "  +0.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.640.41"
$ws.Range("E3").Value = "  +0.32%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "515.44"
$ws.Range("E5").Value = "  +0.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.75"
$ws.Range("E6").Value = "  -0.24%  "

$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.575"
$ws.Range("E8").Value = "  +0.78%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.652.86"
$ws.Range("E9").Value = "  -0.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.30"
$ws.Range("E10").Value = "  -2.49%  "

$ws.Range("E11").Value = "  -0.73%  "

$ws.Range("E12").Value = "  -0.51%  "

$ws.Range("E13").Value = "  +0.86%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.107.93"
$ws.Range("E14").Value = "  +0.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.344.34"
$ws.Range("E15").Value = "  +0.17%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.98"
$ws.Range("E16").Value = "  -0.96%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000137"
$ws.Range("E17").Value = "  -0.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.650.72"
$ws.Range("E18").Value = "  +0.17%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "349.18"
$ws.Range("E19").Value = "  +1.26%  "

$ws.Range("E20").Value = "  -2.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.32"
$ws.Range("E21").Value = "  -1.42%  "

$ws.Range("E22").Value = "  +1.18%  "

$ws.Range("E23").Value = "  -0.36%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.10"
$ws.Range("E24").Value = "  +2.30%  "

$ws.Range("E25").Value = "  -1.59%  "

$ws.Range("E26").Value = "  +3.48%  "

$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0803"
$ws.Range("E28").Value = "  -1.95%  "

$ws.Range("E29").Value = "  -0.96%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.48"
$ws.Range("E31").Value = "  -0.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.92"
$ws.Range("E32").Value = "  -0.25%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.57"
$ws.Range("E33").Value = "  -0.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.86"
$ws.Range("E34").Value = "  +0.49%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.06"
$ws.Range("E35").Value = "  +0.99%  "

$ws.Range("B36").Value = "SuiNetwork"
$ws.Range("C36").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.947"
$ws.Range("E36").Value = "  -11.33%  "

$ws.Range("E37").Value = "  +1.52%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.859"
$ws.Range("E38").Value = "  -0.78%  "

$ws.Range("E39").Value = "  +0.17%  "

$ws.Range("E40").Value = "  +3.99%  "

$ws.Range("E41").Value = "  -1.45%  "

$ws.Range("E42").Value = "  +0.13%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "278.01"
$ws.Range("E43").Value = "  -1.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +0.49%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.602"
$ws.Range("E45").Value = "  -2.45%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.61"
$ws.Range("E46").Value = "  -0.21%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.076.25"
$ws.Range("E47").Value = "  +5.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0530"
$ws.Range("E48").Value = "  -2.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0231"
$ws.Range("E49").Value = "  -0.12%  "

$ws.Range("E50").Value = "  +0.60%  "

$ws.Range("E51").Value = "  -0.57%  "
